# Scheduled runner refresh: re-pull currentAveragePrice / currentAveragePriceNQ /
# currentAveragePriceHQ market data and recompute the dependent Leve price/profit
# columns (H..N) for the affected Leve rows on each job sheet.
# (Cells hold cached numeric snapshots, not live formulas.)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC!row 17 - One for the Road
$ws.Range("H17").Value = 711.8333
$ws.Range("J17").Value = 750.8125
$ws.Range("L17").Value = 2252.4375
$ws.Range("N17").Value = -2588.4375

# ALC!row 28 - The Writing Is Not on the Wall
$ws.Range("H28").Value = 778.10345
$ws.Range("I28").Value = 870.7619
$ws.Range("J28").Value = 534.875
$ws.Range("K28").Value = 870.7619
$ws.Range("L28").Value = 534.875
$ws.Range("M28").Value = -385.7619
$ws.Range("N28").Value = -1504.875

# ALC!row 40 - Stuck in the Moment
$ws.Range("H40").Value = 10091.63
$ws.Range("I40").Value = 1533.091
$ws.Range("J40").Value = 15975.625
$ws.Range("K40").Value = 1533.091
$ws.Range("L40").Value = 15975.625
$ws.Range("M40").Value = -1358.091
$ws.Range("N40").Value = -16325.625

# ALC!row 112 - Making Ends Meet
$ws.Range("H112").Value = 5179.4165
$ws.Range("J112").Value = 5179.4165
$ws.Range("L112").Value = 15538.2495
$ws.Range("N112").Value = -17754.2495

# ALC!row 130 - Technically Still Magic
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# ARM!row 81 - A Halonic Masquerade
$ws.Range("H81").Value = 106737.4
$ws.Range("I81").Value = 81502
$ws.Range("K81").Value = 81502
$ws.Range("M81").Value = -80504

# ARM!row 84 - Why I Wear a Mask (L)
$ws.Range("H84").Value = 106737.4
$ws.Range("I84").Value = 81502
$ws.Range("K84").Value = 244506
$ws.Range("M84").Value = -239514

$ws = $wb.Worksheets.Item("BSM")
# BSM!row 94 - High Steal
$ws.Range("H94").Value = 606.1579
$ws.Range("I94").Value = 601.35297
$ws.Range("K94").Value = 601.35297
$ws.Range("M94").Value = -150.35297

$ws = $wb.Worksheets.Item("CRP")
# CRP!row 58 - You Do the Heavy Lifting
$ws.Range("H58").Value = 4923.909
$ws.Range("I58").Value = 4971.9
$ws.Range("K58").Value = 4971.9
$ws.Range("M58").Value = -4768.9

# CRP!row 86 - Birch, Please
$ws.Range("H86").Value = 11117156
$ws.Range("I86").Value = 16672618
$ws.Range("J86").Value = 6232.3
$ws.Range("K86").Value = 16672618
$ws.Range("L86").Value = 6232.3
$ws.Range("M86").Value = -16671495
$ws.Range("N86").Value = -8478.299999999999

# CRP!row 89 - Built This City on Blocks and Soul (L)
$ws.Range("H89").Value = 11117156
$ws.Range("I89").Value = 16672618
$ws.Range("J89").Value = 6232.3
$ws.Range("K89").Value = 83363090
$ws.Range("L89").Value = 31161.5
$ws.Range("M89").Value = -83357474
$ws.Range("N89").Value = -42393.5

# CRP!row 99 - O Pine
$ws.Range("H99").Value = 3611.4
$ws.Range("I99").Value = 3189.4285
$ws.Range("J99").Value = 4596
$ws.Range("K99").Value = 3189.4285
$ws.Range("L99").Value = 4596
$ws.Range("M99").Value = -1691.4285
$ws.Range("N99").Value = -7592

# CRP!row 126 - A Better Conductor
$ws.Range("H126").Value = 3611.4
$ws.Range("I126").Value = 3189.4285
$ws.Range("J126").Value = 4596
$ws.Range("K126").Value = 9568.2855
$ws.Range("L126").Value = 13788
$ws.Range("M126").Value = -7098.2855
$ws.Range("N126").Value = -18728

# CRP!row 136 - Turali Quality
$ws.Range("H136").Value = 4923.909
$ws.Range("I136").Value = 4971.9
$ws.Range("K136").Value = 14915.7
$ws.Range("M136").Value = -12365.7

$ws = $wb.Worksheets.Item("CUL")
# CUL!row 4 - In Hot Water
$ws.Range("H4").Value = 55444668
$ws.Range("I4").Value = 83474330
$ws.Range("J4").Value = 7393823.5
$ws.Range("K4").Value = 250422990
$ws.Range("L4").Value = 22181470.5
$ws.Range("M4").Value = -250422878
$ws.Range("N4").Value = -22181694.5

# CUL!row 68 - Such a Butter Face
$ws.Range("H68").Value = 2167.1177
$ws.Range("I68").Value = 1843.3334
$ws.Range("J68").Value = 2343.7273
$ws.Range("K68").Value = 5530.0002
$ws.Range("L68").Value = 7031.1819
$ws.Range("M68").Value = -4719.0002
$ws.Range("N68").Value = -8653.1819

# CUL!row 71 - No Margarine of Error (L)
$ws.Range("H71").Value = 2167.1177
$ws.Range("I71").Value = 1843.3334
$ws.Range("J71").Value = 2343.7273
$ws.Range("K71").Value = 16590.0006
$ws.Range("L71").Value = 21093.5457
$ws.Range("M71").Value = -12534.0006
$ws.Range("N71").Value = -29205.5457

# CUL!row 107 - Slippery Service
$ws.Range("H107").Value = 664.8182
$ws.Range("J107").Value = 1399.4
$ws.Range("L107").Value = 4198.200000000001
$ws.Range("N107").Value = -8038.200000000001

# CUL!row 111 - Soup for the Soldier
$ws.Range("H111").Value = 7998
$ws.Range("I111").Value = 6997.25
$ws.Range("K111").Value = 20991.75
$ws.Range("M111").Value = -17924.75

# CUL!row 119 - Super Dark Times
$ws.Range("H119").Value = 10155
$ws.Range("I119").Value = 7592
$ws.Range("K119").Value = 22776
$ws.Range("M119").Value = -17938

$ws = $wb.Worksheets.Item("GSM")
# GSM!row 80 - Needs More Prayerbell
$ws.Range("H80").Value = 1408558.1
$ws.Range("I80").Value = 3091539.5
$ws.Range("J80").Value = 6073.8335
$ws.Range("K80").Value = 3091539.5
$ws.Range("L80").Value = 6073.8335
$ws.Range("M80").Value = -3090541.5
$ws.Range("N80").Value = -8069.8335

# GSM!row 83 - With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 1408558.1
$ws.Range("I83").Value = 3091539.5
$ws.Range("J83").Value = 6073.8335
$ws.Range("K83").Value = 15457697.5
$ws.Range("L83").Value = 30369.1675
$ws.Range("M83").Value = -15452705.5
$ws.Range("N83").Value = -40353.1675

# GSM!row 102 - Put the Metal to the Peddle
$ws.Range("H102").Value = 3576.5
$ws.Range("I102").Value = 3576.5
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 3576.5
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -1954.5
$ws.Range("N102").ClearContents()

# GSM!row 122 - Awarding Academic Excellence
$ws.Range("H122").Value = 5797.4863
$ws.Range("I122").Value = 4175.4165
$ws.Range("J122").Value = 8792.076999999999
$ws.Range("K122").Value = 12526.2495
$ws.Range("L122").Value = 26376.231
$ws.Range("M122").Value = -10076.2495
$ws.Range("N122").Value = -31276.231

# GSM!row 126 - Gold Rush Order
$ws.Range("H126").Value = 5856
$ws.Range("I126").Value = 4798.6
$ws.Range("J126").Value = 8499.5
$ws.Range("K126").Value = 14395.8
$ws.Range("L126").Value = 25498.5
$ws.Range("M126").Value = -11925.8
$ws.Range("N126").Value = -30438.5

$ws = $wb.Worksheets.Item("LTW")
# LTW!row 7 - Tan Before the Ban
$ws.Range("H7").Value = 15676.177
$ws.Range("I7").Value = 14499.777
$ws.Range("J7").Value = 16999.625
$ws.Range("K7").Value = 14499.777
$ws.Range("L7").Value = 16999.625
$ws.Range("M7").Value = -14387.777
$ws.Range("N7").Value = -17223.625

# LTW!row 40 - Best Served Toad
$ws.Range("H40").Value = 8291.959999999999
$ws.Range("I40").Value = 6406.1875
$ws.Range("K40").Value = 6406.1875
$ws.Range("M40").Value = -6270.1875

# LTW!row 122 - Hell on Leather
$ws.Range("H122").Value = 6616.3335
$ws.Range("I122").Value = 6526.1763
$ws.Range("K122").Value = 19578.5289
$ws.Range("M122").Value = -17128.5289

# LTW!row 126 - Battered Books
$ws.Range("H126").Value = 15676.177
$ws.Range("I126").Value = 14499.777
$ws.Range("J126").Value = 16999.625
$ws.Range("K126").Value = 43499.331
$ws.Range("L126").Value = 50998.875
$ws.Range("M126").Value = -41029.331
$ws.Range("N126").Value = -55938.875

$ws = $wb.Worksheets.Item("WVR")
# WVR!row 122 - Heavy Armoire
$ws.Range("H122").Value = 6515.4707
$ws.Range("I122").Value = 5444.875
$ws.Range("K122").Value = 16334.625
$ws.Range("M122").Value = -13884.625
